$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append new log row (row 17) mirroring the "Testmail #2" order/delivery entry.
$ws.Range("A17").Value = "Kun je 10 dozen schroeven bestellen?"
$ws.Range("B17").Value = "mailmind.test@zohomail.eu"
$ws.Range("C17").Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$ws.Range("D17").Value = "Bestelling / Levering"
$ws.Range("E17").Value = "Geachte klant,`nBedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.`nIk zie uw reactie graag tegemoet.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$ws.Range("F17").Value = "2025-06-29 14:43:53"
$ws.Range("G17").Value = "Ja"
$ws.Range("H17").Value = "Ja"
$ws.Range("I17").Value = "Nee"

# Extend the conditional-formatting ranges so the new row gets highlighted too.
$ws.Range("D2:D16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D17"))
$ws.Range("G2:G16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G17"))
$ws.Range("H2:H16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H17"))
$ws.Range("I2:I16").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I17"))

# Update the Dashboard summary count for "Bestelling / Levering".
$dash.Range("B3").Value = 5
